# Update the "Reaction_number" column (C2:C20) on both the NBR and BAR
# sheets with the new sensitivity results.

$wb = $excel.ActiveWorkbook

$nbrValues = @(765, 742, 744, 693, 692, 685, 658, 662, 648, 633, 629, 634, 631, 627, 624, 620, 616, 604, 599)
$barValues = @(694, 680, 669, 657, 657, 655, 666, 659, 674, 670, 676, 664, 663, 659, 649, 645, 645, 646, 642)

$wsNbr = $wb.Worksheets.Item("NBR")
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $wsNbr.Cells.Item($i + 2, 3).Value = $nbrValues[$i]
}

$wsBar = $wb.Worksheets.Item("BAR")
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $wsBar.Cells.Item($i + 2, 3).Value = $barValues[$i]
}
